# "first implementation of Car and movement"
# - Scale the numeric grid on the "data" sheet (B2:E18) by 100 (units: cm -> mm, etc.)
# - Move the active tab / selection from "entry_params" to "data" (cell C3)
# - Move the selection on the "scheme" sheet to L1

$wb = $excel.ActiveWorkbook

# --- data sheet: multiply B2:E18 by 100 ---------------------------------
$wsData = $wb.Worksheets.Item("data")
for ($r = 2; $r -le 18; $r++) {
    for ($c = 2; $c -le 5; $c++) {
        $cell = $wsData.Cells.Item($r, $c)
        $v = $cell.Value2
        if ($v -ne $null) {
            $cell.Value = $v * 100
        }
    }
}

# --- scheme sheet: move selection to L1 ---------------------------------
$wsScheme = $wb.Worksheets.Item("scheme")
$wsScheme.Activate()
$wsScheme.Range("L1").Select()

# --- data sheet: becomes the active tab, selection at C3 ----------------
$wsData.Activate()
$wsData.Range("C3").Select()

Write-Output "done"
